# Update cryptos list values (Price and Volume(1h) columns) to match the latest fetch.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.609.45"
$ws.Range("E2").Value = "  +1.63%  "

$ws.Range("D3").Value = "1.826.97"
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.23"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5308"
$ws.Range("E7").Value = "  -2.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3970"
$ws.Range("E8").Value = "  +4.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07758"
$ws.Range("E9").Value = "  +3.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.02"
$ws.Range("E10").Value = "  -0.05%  "

$ws.Range("E11").Value = "  +2.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.15"
$ws.Range("E12").Value = "  +2.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.320"
$ws.Range("E13").Value = "  +1.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.592"
$ws.Range("E14").Value = "  +3.07%  "

$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").Value = "1.827.63"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.06"
$ws.Range("E17").Value = "  +3.43%  "

$ws.Range("E18").Value = "  +2.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06603"
$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.81"
$ws.Range("E20").Value = "  +1.71%  "

$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("E22").Value = "  +2.65%  "

$ws.Range("D23").Value = "28.618.04"
$ws.Range("E23").Value = "  +1.50%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.235"
$ws.Range("E25").Value = "  +6.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.75"
$ws.Range("E26").Value = "  +1.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.98"
$ws.Range("E27").Value = "  +0.73%  "

$ws.Range("D28").Value = "2.035.62"
$ws.Range("E28").Value = "  +1.17%  "

$ws.Range("E29").Value = "  +3.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.69"
$ws.Range("E30").Value = "  +2.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.149"
$ws.Range("E31").Value = "  +2.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1127"
$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.737"
$ws.Range("E33").Value = "  +2.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.653"
$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07320"
$ws.Range("E35").Value = "  +4.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2265"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02351"
$ws.Range("E37").Value = "  +1.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.898"
$ws.Range("E38").Value = "  +4.88%  "

$ws.Range("E39").Value = "  +2.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.39"
$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6291"
$ws.Range("E41").Value = "  +1.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.194"
$ws.Range("E42").Value = "  +1.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.397"
$ws.Range("E44").Value = "  -1.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.53"
$ws.Range("E45").Value = "  +1.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5933"
$ws.Range("E46").Value = "  +2.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.721"
$ws.Range("E47").Value = "  +1.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.44"
$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.998"
$ws.Range("E49").Value = "  +3.79%  "

$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("E51").Value = "  +1.85%  "

